$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change: remove "/RME" from the steel (S) description line ---
$newText = "10% CR/LFM+CDN/H:2`n25% CR+PC/LFM+CDN/H:1`n30% S/LFM+CDN/H:1`n8% S+SL/LFM+CDN/H:1`n10% S/LFBR+CDN/H:1`n15% W/LWAL+CDN/H:1`n2% MUR/LWAL+CDN/H:1"
$ws.Range("B2").Value = $newText

# --- Formatting: wrap the long multi-line text and size the row to fit it ---
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 365

# --- View state: select B2:B12 (active cell ends up at B12, as in the saved file) ---
$ws.Range("B12").Select() | Out-Null
$ws.Range("B2:B12").Select() | Out-Null
